$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric value updates ---
$ws.Range("B2").Value = 23589205.05
$ws.Range("B3").Value = 12798874.58
$ws.Range("B4").Value = 10790330.47
$ws.Range("B7").Value = 226329478.62
$ws.Range("B8").Value = 106835116.29
$ws.Range("B9").Value = 119494362.33
$ws.Range("B13").Value = 82930332.0071497
$ws.Range("B20").Value = 97955830.00999998
$ws.Range("B21").Value = 16500000
$ws.Range("B24").Value = 38.36175787029699
$ws.Range("B25").Value = 54.5910283895756
$ws.Range("B26").Value = 23.85181807547456
$ws.Range("B49").Value = 28.56709602348285
$ws.Range("B52").Value = 730.4072408515789
$ws.Range("B55").Value = 2260.582755498804

# --- Row 34: "Operational Liquidity / Actuals" row loses its numeric
#     value and A34 picks up the bold header style used elsewhere
#     (copy formats from A33, which already carries that style). ---
$ws.Range("B34").ClearContents()
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)

# --- Percentage cells stored as literal text, e.g. "29.67%".
#     Typing a percent-looking string into a cell makes Excel coerce it
#     to a formatted number, so force Text format first, then strip the
#     resulting formatting back off (ClearFormats) so the cell keeps the
#     literal string value with the original (default) style. ---
function Set-PercentText($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

Set-PercentText "B38" "29.67%"
Set-PercentText "B42" "24.90%"
Set-PercentText "B44" "3.26%"
Set-PercentText "B59" "5.12%"
Set-PercentText "B62" "11.90%"
Set-PercentText "B65" "3.73%"
Set-PercentText "B68" "4.33%"
Set-PercentText "B71" "12.89%"
Set-PercentText "B74" "3.20%"
Set-PercentText "B86" "0.03%"
Set-PercentText "B176" "27.61%"
Set-PercentText "B215" "11.25%"
Set-PercentText "B222" "20.23%"
Set-PercentText "B226" "21.80%"
Set-PercentText "B228" "33.20%"
Set-PercentText "B230" "22.01%"
Set-PercentText "B234" "21.24%"
Set-PercentText "B238" "20.20%"
Set-PercentText "B242" "19.50%"
Set-PercentText "B243" "27.20%"
Set-PercentText "B244" "31.80%"
Set-PercentText "B246" "18.54%"
Set-PercentText "B247" "26.10%"
Set-PercentText "B248" "30.50%"
